$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.532.47"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "3.378.43"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "404.84"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.38"
$ws.Range("E6").Value = "  +7.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.667"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.119"
$ws.Range("E10").Value = "  -7.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.35"
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "3.910.54"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.40"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.74"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "3.384.54"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "61.558.36"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.01"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.95"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("E20").Value = "  -7.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.18"
$ws.Range("E21").Value = "  -4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "84.94"
$ws.Range("E22").Value = "  +4.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "316.50"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.74"
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.11"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  +11.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.44"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.24"
$ws.Range("E28").Value = "  +5.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.62"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.36"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.38"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0478"
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.69"
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.40"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "139.14"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.293"
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.97"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.63"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.22"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").Value = "2.121.75"
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("E50").Value = "  -7.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.88"
$ws.Range("E51").Value = "  +0.80%  "
